$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Added connection between cost and demand" / variable+fixed cost inputs:
# the container_height input (B3) feeds the fluidMass formula in B6
# (=(container_height*3.14*(container_diameter/2)*B5)/1000), so updating
# B3 automatically ripples into the cached formula result in B6.
$ws.Range("B3").Value = 24.999980000000001
